$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "proklicovani" (sheet2): add the "laborka"/"venku" control rows and
# their combined totals underneath the existing kontrola/cm1/cm3/cm5 block.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("proklicovani")

# Seed the new shared-string table entries in the exact order the original
# author's workbook uses them (1cm, 3cm, 5cm, laborka, venku, p, ...) so the
# rebuilt sharedStrings.xml lines up index-for-index.
$ws2.Range("H26").Value = "1cm"
$ws2.Range("H27").Value = "3cm"
$ws2.Range("H28").Value = "5cm"
$ws2.Range("F19").Value = "laborka"
$ws2.Range("F25").Value = "venku"
$ws2.Range("L18").Value = "p"

# New "venku" (outdoor) block, rows 25-28, mirrors rows 19-22.
$ws2.Range("H25").Value = "kontrola"
$ws2.Range("I25").Value = 40
$ws2.Range("J25").Value = 35
$ws2.Range("K25").Value = 75
$ws2.Range("L25").Formula = "=I25/K25"

$ws2.Range("I26").Value = 26
$ws2.Range("J26").Value = 49
$ws2.Range("K26").Value = 75
$ws2.Range("L26:L28").Formula = "=I26/K26"

$ws2.Range("I27").Value = 0
$ws2.Range("J27").Value = 75
$ws2.Range("K27").Value = 75

$ws2.Range("I28").Value = 0
$ws2.Range("J28").Value = 75
$ws2.Range("K28").Value = 75

# Combined totals (laborka + venku), rows 32-35.
$ws2.Range("H32").Value = "kontrola"
$ws2.Range("I32").Formula = "=I19+I25"
$ws2.Range("J32:K32").Formula = "=J19+J25"
$ws2.Range("L32").Formula = "=I32/K32"

$ws2.Range("H33").Value = "1cm"
$ws2.Range("I33:K35").Formula = "=I20+I26"
$ws2.Range("L33:L35").Formula = "=I33/K33"

$ws2.Range("H34").Value = "3cm"

$ws2.Range("H35").Value = "5cm"

# ---------------------------------------------------------------------------
# New worksheet "List1" (sheet3) with the raw germination data.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "List1"
$ws3.PageSetup.LeftMargin = 0.7 * 72
$ws3.PageSetup.RightMargin = 0.7 * 72
$ws3.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws3.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws3.PageSetup.HeaderMargin = 0.3 * 72
$ws3.PageSetup.FooterMargin = 0.3 * 72

$ws3.Range("A1").Value = "germ"
$ws3.Range("B1").Value = "n"
$ws3.Range("C1").Value = "exp"
$ws3.Range("D1").Value = "treatment"

$ws3.Range("A2").Value = 10
$ws3.Range("B2").Value = 75
$ws3.Range("C2").Value = "lab"
$ws3.Range("D2").Value = 0

$ws3.Range("A3").Value = 7
$ws3.Range("B3").Value = 75
$ws3.Range("C3").Value = "lab"
$ws3.Range("D3").Value = 1

$ws3.Range("A4").Value = 0
$ws3.Range("B4").Value = 75
$ws3.Range("C4").Value = "lab"
$ws3.Range("D4").Value = 3

$ws3.Range("A5").Value = 0
$ws3.Range("B5").Value = 75
$ws3.Range("C5").Value = "lab"
$ws3.Range("D5").Value = 5

$ws3.Range("A6").Value = 40
$ws3.Range("B6").Value = 75
$ws3.Range("C6").Value = "ext"
$ws3.Range("D6").Value = 0

$ws3.Range("A7").Value = 26
$ws3.Range("B7").Value = 75
$ws3.Range("C7").Value = "ext"
$ws3.Range("D7").Value = 1

$ws3.Range("A8").Value = 0
$ws3.Range("B8").Value = 75
$ws3.Range("C8").Value = "ext"
$ws3.Range("D8").Value = 3

$ws3.Range("A9").Value = 0
$ws3.Range("B9").Value = 75
$ws3.Range("C9").Value = "ext"
$ws3.Range("D9").Value = 5

# Match the author's final on-screen state: "proklicovani" keeps F19:L35
# selected, and the newly added "List1" ends up as the active tab.
$ws2.Activate() | Out-Null
$ws2.Range("F19:L35").Select() | Out-Null
$ws3.Activate() | Out-Null
